# Update the R6 vs Legacy comparison table on the "RLcomp_valid" sheet so
# the three R6-syntax examples that contain a literal "$" escape it the
# same way every other entry in the table already does (e.g. PM_data\$new()).
#
# Order matters: new shared-string entries are appended in first-write
# order, so we touch B6, then B5, then B3 to reproduce the exact shared
# string table ordering of the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RLcomp_valid")

$ws.Range("B6").Value = "PM_result\`$step()"
$ws.Range("B5").Value = "PM_valid\`$plot()"
$ws.Range("B3").Value = "PM_result\`$op\`$plot(resid = T,…)"

# Restore a single-cell selection on B4 (was previously B4:C4).
$ws.Activate()
$ws.Range("B4").Select()
